# Apply updates described in the commit "Minor changes to
# data_preprocessing_into_inputs.py" to the RESERVE_settings workbook.

$wb = $excel.ActiveWorkbook

$wsReserve = $wb.Worksheets.Item("RESERVE Settings")
$wsInput   = $wb.Worksheets.Item("Input Data Settings")

# --- Input Data Settings sheet ----------------------------------------------
# Load forecast data source file
$wsInput.Range("A3").Value = "PSE_RTPD_load_forecast.csv"

# Load actual data source file
$wsInput.Range("A2").Value = "PSE_RTD_load_forecast.csv"

# Wind actual data source file
$wsInput.Range("A4").Value = "PSE_wind_5_minute_actuals.csv"

# Solar actual data source file
$wsInput.Range("A6").Value = "PSE_solar_5_minute_actuals.csv"

# --- RESERVE Settings sheet -------------------------------------------------
# ANCHOR_DATE value: 2017-01-01 -> 2019-01-01
$wsReserve.Range("C3").Value = [DateTime]"2019-01-01"

# LATITUDE value: 36.6777 -> 47.544
$wsReserve.Range("C6").Value = 47.544

# LONGITUDE value: -119.4179 -> -120.411
$wsReserve.Range("C7").Value = -120.411

# MODEL_NAME value: rescue_v1_4_manually_cleaned -> reserve_PSE_2022
$wsReserve.Range("C2").Value = "reserve_PSE_2022"

# --- Window / selection state ------------------------------------------------
# Active tab moves from "Input Data Settings" back to "RESERVE Settings";
# selections are updated to reflect where the author last clicked.
[void]$wsInput.Range("C16").Select()
[void]$wsReserve.Activate()
[void]$wsReserve.Range("C2").Select()
